$wb = $excel.ActiveWorkbook
$sheetA = $wb.Worksheets.Item("a")
$newSheet = $wb.Worksheets.Add($null, $sheetA)
$newSheet.Name = "b"
$sheetA.Range("A1:AD12").Copy($newSheet.Range("A1"))
$newSheet.Columns("E").Delete()
$newSheet.Columns("C").Delete()
$newSheet.Range("AA3:AA12").ClearContents()
$newSheet.Range("AB6").ClearContents()

# Try toggling interior color on/off to nudge applyFill flag
$newSheet.Range("AB1").Interior.ColorIndex = 3
$newSheet.Range("AB1").Interior.ColorIndex = -4142
$newSheet.Range("AA2").Interior.ColorIndex = 3
$newSheet.Range("AA2").Interior.ColorIndex = -4142
$newSheet.Range("AB2").Interior.ColorIndex = 3
$newSheet.Range("AB2").Interior.ColorIndex = -4142

Write-Output "done"
